$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Remove SK2 from products": the SK2 (old data) product lived in column G.
# Delete it, which shifts the remaining product columns (H:K) left by one.
$ws.Columns("G:G").Delete()

# The "Standard Kalk Kat3" product column (originally K, now J after the
# delete above) is moved to be the first product column, right after
# Description, i.e. column C.
$ws.Columns("J:J").Cut()
$ws.Columns("C:C").Insert()

# Restore the selection to match the saved view state.
$ws.Range("E23").Select() | Out-Null
